$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-02 Saturday", "2024-03-03 Sunday"),
    @("112÷8=", "561÷2="),
    @("428÷9=", "828÷2="),
    @("289÷8=", "292÷5="),
    @("262÷4=", "751÷7="),
    @("110÷6=", "799÷2="),
    @("602÷7=", "572÷4="),
    @("706÷3=", "733÷9="),
    @("106÷2=", "540÷2="),
    @("785÷4=", "612÷9="),
    @("995÷2=", "566÷6="),
    @("113÷8=", "151÷2="),
    @("234÷8=", "356÷3="),
    @("429÷5=", "951÷8="),
    @("224÷5=", "241÷5="),
    @("266÷7=", "212÷7="),
    @("851÷9=", "501÷6="),
    @("556÷7=", "474÷3="),
    @("154÷4=", "586÷3="),
    @("181÷2=", "600÷9="),
    @("811÷4=", "305÷6="),
    @("197÷6=", "346÷3="),
    @("226÷2=", "933÷4="),
    @("421÷5=", "180÷2="),
    @("378÷5=", "180÷6="),
    @("118÷7=", "662÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
